$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header renames (A1 stays "Unnamed: 1")
$ws.Range("B1").Value = "Unnamed: 4"
$ws.Range("C1").Value = "Unnamed: 16"
$ws.Range("D1").Value = "Unnamed: 19"

# Row 2: was blank row (B2 already blank/text -> leave), now PREDICTION / Nr of points / Points
$ws.Range("A2").Value = "PREDICTION"
$ws.Range("C2").Value = "Nr of points"
$ws.Range("D2").Value = "Points"

# Row 3: was "Predicted headform score..." label row -> now D Green data row
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = "D Green"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Row 4: was fully blank already (A4/B4/C4/D4 all blank text) -> now Green data row
$ws.Range("B4").Value = "Green"
$ws.Range("C4").Value = 29
$ws.Range("D4").Value = 29

# Row 5: was "VERIFICATION" label row -> now Yellow data row
$ws.Range("A5").Value = $null
$ws.Range("B5").Value = "Yellow"
$ws.Range("C5").Value = 101
$ws.Range("D5").Value = 75.75

# Row 6: was "Testpoint"/"Prediction" header row -> now Orange data row
$ws.Range("A6").Value = $null
$ws.Range("B6").Value = "Orange"
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 22

# Row 7: was "7,2"/"Green" -> now Brown data row
$ws.Range("A7").Value = $null
$ws.Range("B7").Value = "Brown"
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 5

# Row 8: was "13,-2"/"Yellow" -> now Red data row
$ws.Range("A8").Value = $null
$ws.Range("B8").Value = "Red"
$ws.Range("C8").Value = 26
$ws.Range("D8").Value = 0

# Row 9: was "9,1"/"Orange" -> now Default Red data row
$ws.Range("A9").Value = $null
$ws.Range("B9").Value = "Default Red"
$ws.Range("C9").Value = 12
$ws.Range("D9").Value = 0

# Row 10: was "11,7"/"Red" -> now Blue data row (D10 stays blank/text)
$ws.Range("A10").Value = $null
$ws.Range("B10").Value = "Blue"
$ws.Range("C10").Value = 0

# Row 11: was "2,6"/"Yellow" -> now the "Predicted headform score..." label row (B11 stays blank/text)
$ws.Range("A11").Value = "Predicted headform score (excluding blue points)"
$ws.Range("C11").Value = 232
$ws.Range("D11").Value = 131.75

# Row 12 (new row): VERIFICATION label
$ws.Range("A12").Value = "VERIFICATION"

# Column E no longer exists in the new layout (A1:D12) -> delete it
$ws.Range("E1:E11").Delete()
